# Emissionskatalog.xlsx — add the "reference" URL to every data row (2-200)
# in column G of the "Tabelle1" sheet (commit: "All Emissions from
# Emissionskatalog.xlsx in DB").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.zwei-grad-eine-tonne.at/hintergrund-berechnungen/abschnitt-i-lustvoll-die-welt-retten"

$lastRow = 200

$rng = $ws.Range("G2:G$lastRow")
$rng.Value = $url
